# Fruta / hortaliza, semanal
# Insert 6 new weekly price rows for Limón (Vega Central Mapocho de Santiago)
# right before the current row 1158, pushing the existing data down by 6 rows
# (old row N becomes row N+6), and populate the newly-inserted rows with the
# new week's data (date serial 44516 = 2021-11-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 1158..1230 down to 1164..1236 by inserting 6 blank rows at 1158.
$ws.Rows("1158:1163").Insert()

# Common (unchanged) columns for every row in this dataset.
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102003
$categoria = "Limón"

$newRows = @(
    @{ Row=1158; Fecha=44516; Variedad="Sin especificar"; Calidad="1a amarillo"; Volumen=440; PrecioMin=7000; PrecioMax=7000; PrecioProm=7000; Unidad="$/malla 18 kilos"; Origen="Provincia del Elquí"; PrecioKg=389; KgUnidad=18 },
    @{ Row=1159; Fecha=44516; Variedad="Sin especificar"; Calidad="1a amarillo"; Volumen=530; PrecioMin=6500; PrecioMax=7000; PrecioProm=6764; Unidad="$/malla 18 kilos"; Origen="Región de O'Higgins"; PrecioKg=376; KgUnidad=18 },
    @{ Row=1160; Fecha=44516; Variedad="Sin especificar"; Calidad="2a amarillo"; Volumen=350; PrecioMin=6000; PrecioMax=6000; PrecioProm=6000; Unidad="$/malla 18 kilos"; Origen="Provincia del Elquí"; PrecioKg=333; KgUnidad=18 },
    @{ Row=1161; Fecha=44516; Variedad="Sin especificar"; Calidad="2a amarillo"; Volumen=520; PrecioMin=5500; PrecioMax=6000; PrecioProm=5788; Unidad="$/malla 18 kilos"; Origen="Región de O'Higgins"; PrecioKg=322; KgUnidad=18 },
    @{ Row=1162; Fecha=44516; Variedad="Sin especificar"; Calidad="3a amarillo"; Volumen=300; PrecioMin=4500; PrecioMax=4500; PrecioProm=4500; Unidad="$/malla 18 kilos"; Origen="Provincia del Elquí"; PrecioKg=250; KgUnidad=18 },
    @{ Row=1163; Fecha=44516; Variedad="Sin especificar"; Calidad="3a amarillo"; Volumen=600; PrecioMin=4000; PrecioMax=4500; PrecioProm=4250; Unidad="$/malla 18 kilos"; Origen="Región de O'Higgins"; PrecioKg=236; KgUnidad=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PrecioMin
    $ws.Cells.Item($row, 15).Value = $r.PrecioMax
    $ws.Cells.Item($row, 16).Value = $r.PrecioProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
